$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.197.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.354.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.61%  "

$ws.Range("E13").Value = "  +3.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.724.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.379.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.795"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.186.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("E20").Value = "  +4.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.15"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0725"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -36.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.946.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "

$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("E45").Value = "  +4.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.90%  "

$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.586.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.87"
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = "  -4.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
